$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '98.463.29'
$ws.Range("E2").Value = '  +0.00%  '

# Row 3
$ws.Range("D3").Value = '3.354.15'
$ws.Range("E3").Value = '  +0.61%  '

# Row 4
$ws.Range("E4").Value = '  -0.02%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '257.00'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.46%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '664.15'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +6.05%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.52'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +8.60%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.474'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +22.30%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '1.07'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +23.74%  '

# Row 10
$ws.Range("E10").Value = '  -0.02%  '

# Row 11
$ws.Range("D11").Value = '3.347.26'
$ws.Range("E11").Value = '  +0.50%  '

# Row 12
$ws.Range("E12").Value = '  +8.40%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '42.26'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +13.83%  '

# Row 14
$ws.Range("E14").Value = '  +10.28%  '

# Row 15
$ws.Range("D15").Value = '98.418.75'
$ws.Range("E15").Value = '  +0.24%  '

# Row 16
$ws.Range("B16").Value = 'Toncoin'
$ws.Range("C16").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '5.69'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +3.42%  '

# Row 17
$ws.Range("B17").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C17").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D17").Value = '3.978.17'
$ws.Range("E17").Value = '  -0.04%  '

# Row 18
$ws.Range("D18").Value = '3.356.83'
$ws.Range("E18").Value = '  +0.51%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '7.65'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +26.30%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '16.70'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +10.36%  '

# Row 21
$ws.Range("B21").Value = 'BitcoinCash'
$ws.Range("C21").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '529.87'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +8.51%  '

# Row 22
$ws.Range("B22").Value = 'SuiNetwork'
$ws.Range("C22").Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '3.57'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.91%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '10.59'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +12.87%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.0000218'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +3.34%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.434'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +52.29%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '102.08'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +15.13%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '6.17'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +9.77%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '12.50'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +5.12%  '

# Row 29
$ws.Range("D29").Value = '3.539.11'
$ws.Range("E29").Value = '  +0.73%  '

# Row 30
$ws.Range("E30").Value = '  +7.47%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '10.99'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +14.25%  '

# Row 33
$ws.Range("E33").Value = '  -1.60%  '

# Row 34
$ws.Range("E34").Value = '  +0.13%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '29.29'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +5.33%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.537'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +16.94%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '7.79'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +6.79%  '

# Row 38
$ws.Range("E38").Value = '  +8.59%  '

# Row 39
$ws.Range("E39").Value = '  +5.34%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '524.32'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +5.61%  '

# Row 41
$ws.Range("B41").Value = 'WhiteBITCoin'
$ws.Range("C41").Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '24.70'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -0.59%  '

# Row 42
$ws.Range("B42").Value = 'Fetch.AI'
$ws.Range("C42").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.33'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +5.43%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '3.85'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +4.03%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0432'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +32.94%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '3.41'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +3.50%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.821'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +5.22%  '

# Row 47
$ws.Range("E47").Value = '  -0.01%  '

# Row 48
$ws.Range("B48").Value = 'Filecoin'
$ws.Range("C48").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '5.14'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +11.63%  '

# Row 49
$ws.Range("B49").Value = 'Stacks'
$ws.Range("C49").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.05'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +6.48%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '7.85'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +17.33%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '50.95'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +11.84%  '
